# BCF-API_URLs.docx edit
#
# Renames the "(Topics) Header" endpoint section to "(Topics) Revisions":
#   - the section's title row text
#   - the two example request URLs ending in "}/header" (one with, one
#     without, a leading space before the closing brace) become
#     ".../revisions"
#
# Plain text Find & Replace is used throughout so that the existing run
# formatting (bold/italic/accent colour on the title, "en-US" language
# tag on the URL runs, the xml:space="preserve" leading-space run, etc.)
# is left completely untouched - only the literal characters change.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Section title: "(Topics) Header" -> "(Topics) Revisions"
#    There is exactly one "(Topics) Header" in the document, so a plain
#    Replace is unambiguous. MatchCase keeps this from interacting with
#    any other casing of the word.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "(Topics) Header", $true, $false, $false, $false, $false,
    $true, 1, $false, "(Topics) Revisions", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) The example URLs: "...}/header" -> "...}/revisions".
#    One call with Replace = wdReplaceAll (2) catches both the
#    no-leading-space and the leading-space variant in one pass (the
#    no-leading-space pattern is a substring of the leading-space one),
#    without touching the unrelated, pre-existing ".../revisions"
#    endpoints used elsewhere in the document (their text never
#    contained "header").
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "}/header", $true, $false, $false, $false, $false,
    $true, 1, $false, "}/revisions", 2) | Out-Null
